$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 2962
$ws.Range("I2").Value = 8058
$ws.Range("J2").Value = 33047
$ws.Range("K2").Value = 206
$ws.Range("L2").Value = 9167
$ws.Range("M2").Value = 601
$ws.Range("N2").Value = 5680
$ws.Range("O2").Value = 21
$ws.Range("P2").Value = 124
$ws.Range("Q2").Value = 53
$ws.Range("R2").Value = 438
$ws.Range("S2").Value = 3578
$ws.Range("T2").Value = 5815
$ws.Range("U2").Value = 418
$ws.Range("V2").Value = 51385
$ws.Range("W2").Value = 12
$ws.Range("X2").Value = 51184
$ws.Range("Y2").Value = 79
$ws.Range("Z2").Value = 795
$ws.Range("AA2").Value = 357
